# Update the BOM entry for Y1 (Crystal or Oscillator): replace the old
# manufacturer part number "830207390509" with the new one
# "ASEKDV-32.768kHz-LC-T". This value is shared across the Comment (A),
# DesignItemId (F) and Footprint (G) columns for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "830207390509"
$newValue = "ASEKDV-32.768kHz-LC-T"

# Use a cell in the same row that is NOT being modified as a format
# donor, so that after changing the text we can restore the original
# cell formatting (this runtime does not preserve the "number stored as
# text" quote-prefix style when a text value is written via Value2).
$formatDonor = $ws.Range("B18")

$usedRange = $ws.UsedRange
foreach ($cell in $usedRange.Cells) {
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue

        $formatDonor.Copy()
        $cell.PasteSpecial(-4122)  # xlPasteFormats
    }
}

$excel.CutCopyMode = 0
